$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "43.169.52"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +1.38%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.274.88"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +1.64%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.01"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.24%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "113.72"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.22%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "303.63"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +7.32%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.630"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +1.03%  "

$ws.Range("E8").Value = "  -0.21%  "

$ws.Range("E9").Value = "  +0.47%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "44.70"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -4.40%  "

$ws.Range("E11").Value = "  -0.01%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "54.93"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +1.40%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "8.97"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -2.13%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "1.06"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +19.98%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.104"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.71%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "15.54"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +1.24%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "2.613.45"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +1.47%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.304.90"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +3.04%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "43.160.68"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.90%  "

$ws.Range("E20").Value = "  +0.34%  "

$ws.Range("E21").Value = "  +5.83%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "75.22"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +4.14%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "3.56"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +12.14%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "258.31"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +11.35%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.46"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +4.57%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "9.07"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -2.46%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "11.71"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -3.11%  "

$ws.Range("E28").Value = "  -0.18%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.24"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.23%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "38.29"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -5.26%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "22.37"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +5.76%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "175.59"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.61%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.18"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -3.20%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.0899"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.17%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "5.74"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +2.72%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "5.12"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +9.86%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "4.29"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -5.80%  "

$ws.Range("E38").Value = "  +0.55%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.0379"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +1.88%  "

$ws.Range("E40").Value = "  -1.16%  "

$ws.Range("E41").Value = "  -5.31%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "72.46"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.07%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.234"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.07%  "

$ws.Range("E44").Value = "  -0.09%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "12.68"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -5.93%  "

$ws.Range("E46").Value = "  +0.97%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "5.61"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.85%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "107.92"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +7.07%  "

$ws.Range("E49").Value = "  +1.86%  "

$ws.Range("E50").Value = "  +2.39%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "73.69"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +5.15%  "
